# Update cfb_weather.xlsx with Timestamp 2024-12-06T10:01:09.971485
# Refreshes weather/wind readings for both sheets and re-orders the
# "Other" sheet's games to match the latest scrape ordering.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FBS")
$ws2 = $wb.Worksheets.Item("Other")

# --- Sheet FBS: refresh weather data + timestamp (rows 2-7) ---
# Row 2
$ws1.Range("N2").Value = "ENE"
$ws1.Range("O2").Value = 29.06
$ws1.Range("P2").Value = 7.2
$ws1.Range("Q2").Value = "ENE"
$ws1.Range("S2").Value = -0.12
$ws1.Range("T2").Value = -0.37
$ws1.Range("U2").Value = 1.4
$ws1.Range("V2").Value = "41.3874924, -73.9640891"
$ws1.Range("AK2").Value = "2024-12-06T10:01:09.971485"
# Row 3
$ws1.Range("N3").Value = "SE"
$ws1.Range("O3").Value = 29.36
$ws1.Range("P3").Value = 2.7
$ws1.Range("Q3").Value = "SSE"
$ws1.Range("S3").Value = -0.08
$ws1.Range("T3").Value = -0.33
$ws1.Range("U3").Value = -4.1
$ws1.Range("V3").Value = "43.6028839, -116.1958882"
$ws1.Range("AK3").Value = "2024-12-06T10:01:09.971485"
# Row 4
$ws1.Range("M4").Value = "SSE"
$ws1.Range("N4").Value = "SSE"
$ws1.Range("P4").Value = 4.8
$ws1.Range("Q4").Value = "SSE"
$ws1.Range("U4").Value = 0
$ws1.Range("V4").Value = "33.8201052, -85.76647"
$ws1.Range("Z4").Value = -110
$ws1.Range("AK4").Value = "2024-12-06T10:01:09.971485"
# Row 5
$ws1.Range("M5").Value = "NE"
$ws1.Range("N5").Value = "NE"
$ws1.Range("O5").Value = 36.8
$ws1.Range("P5").Value = 11.4
$ws1.Range("Q5").Value = "NE"
$ws1.Range("U5").Value = 1.1
$ws1.Range("V5").Value = "39.5197009, -84.7330255"
$ws1.Range("AK5").Value = "2024-12-06T10:01:09.971485"
# Row 6
$ws1.Range("M6").Value = "SW"
$ws1.Range("N6").Value = "SW"
$ws1.Range("O6").Value = 47.72
$ws1.Range("P6").Value = 4
$ws1.Range("Q6").Value = "WSW"
$ws1.Range("U6").Value = -5.9
$ws1.Range("V6").Value = "30.2158434, -92.0417371"
$ws1.Range("AK6").Value = "2024-12-06T10:01:09.971485"
# Row 7
$ws1.Range("M7").Value = "N"
$ws1.Range("N7").Value = "N"
$ws1.Range("O7").Value = 46.76
$ws1.Range("P7").Value = 8.1
$ws1.Range("Q7").Value = "N"
$ws1.Range("R7").Value = 0.8
$ws1.Range("U7").Value = -2.2
$ws1.Range("V7").Value = "32.8377223, -96.7827859"
$ws1.Range("AK7").Value = "2024-12-06T10:01:09.971485"

# --- Sheet Other: reorder games + refresh weather data (rows 2-5) ---
# Row 2
$ws2.Range("A2").Value = "Montana vs South Dakota State"
$ws2.Range("B2").Value = "South Dakota State"
$ws2.Range("C2").Value = "Montana"
$ws2.Range("D2").Value = "SAT 12/07"
$ws2.Range("E2").Value = "01:00 PM"
$ws2.Range("F2").Value = "Mid"
$ws2.Range("J2").Value = -474.5684815
$ws2.Range("K2").Value = 46.7
$ws2.Range("L2").Value = 47.64
$ws2.Range("N2").Value = 2016
$ws2.Range("O2").Value = "E"
$ws2.Range("P2").Value = "E"
$ws2.Range("Q2").Value = 47.36000000000001
$ws2.Range("R2").Value = 11.4
$ws2.Range("S2").Value = "E"
$ws2.Range("T2").Value = 0
$ws2.Range("U2").Value = 0
$ws2.Range("V2").Value = 0
$ws2.Range("X2").Value = "44.3210182, -96.7801386"
# Row 3
$ws2.Range("A3").Value = "Rhode Island vs Mercer"
$ws2.Range("B3").Value = "Mercer"
$ws2.Range("C3").Value = "Rhode Island"
$ws2.Range("D3").Value = "SAT 12/07"
$ws2.Range("E3").Value = "02:00 PM"
$ws2.Range("F3").Value = "Low"
$ws2.Range("J3").Value = 105.98195272
$ws2.Range("K3").Value = 64.83
$ws2.Range("L3").Value = 52.81
$ws2.Range("N3").Value = 2013
$ws2.Range("O3").Value = "ENE"
$ws2.Range("P3").Value = "ENE"
$ws2.Range("Q3").Value = 56.12
$ws2.Range("R3").Value = 5.3
$ws2.Range("S3").Value = "ENE"
$ws2.Range("T3").Value = 0
$ws2.Range("U3").Value = 0
$ws2.Range("V3").Value = 0
$ws2.Range("X3").Value = "32.8262075, -83.6522485"
# Row 4
$ws2.Range("A4").Value = "Villanova vs Incarnate Word"
$ws2.Range("B4").Value = "Incarnate Word"
$ws2.Range("C4").Value = "Villanova"
$ws2.Range("D4").Value = "SAT 12/07"
$ws2.Range("E4").Value = "01:00 PM"
$ws2.Range("F4").Value = "Low"
$ws2.Range("J4").Value = 81.08228299999999
$ws2.Range("K4").Value = 70.74
$ws2.Range("L4").Value = 55.05
$ws2.Range("N4").Value = 2008
$ws2.Range("O4").Value = "SW"
$ws2.Range("P4").Value = "SW"
$ws2.Range("Q4").Value = 48.38
$ws2.Range("R4").Value = 5.7
$ws2.Range("S4").Value = "SW"
$ws2.Range("T4").Value = 0
$ws2.Range("U4").Value = 0
$ws2.Range("V4").Value = 0
$ws2.Range("X4").Value = "29.4674787, -98.470014"
# Row 5
$ws2.Range("A5").Value = "Illinois State vs UC Davis"
$ws2.Range("B5").Value = "UC Davis"
$ws2.Range("C5").Value = "Illinois State"
$ws2.Range("D5").Value = "SAT 12/07"
$ws2.Range("E5").Value = "01:00 PM"
$ws2.Range("F5").Value = "High"
$ws2.Range("J5").Value = -231.4896765
$ws2.Range("K5").Value = 62.21
$ws2.Range("L5").Value = 53.1
$ws2.Range("N5").Value = 2007
$ws2.Range("O5").Value = "N"
$ws2.Range("P5").Value = "NNW"
$ws2.Range("Q5").Value = 61.1
$ws2.Range("R5").Value = 1.2
$ws2.Range("S5").Value = "N"
$ws2.Range("T5").Value = 0
$ws2.Range("U5").Value = 0
$ws2.Range("V5").Value = 0
$ws2.Range("X5").Value = "38.5365266, -121.7627936"
